# CLEAN PROJECT FOLDER !
# Update the Sheet1 task tracker: close out several tasks to "ok", flip the
# "DELAY" status to "Tạm ổn", mark one task as "nearly done" and append three
# new module tasks (rows 13-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Phase 1: copy existing cell formatting onto the cells whose status
# color needs to change, using the still-unmodified template cells
# (G7/G11 currently hold the "Tạm ổn" / "ok" looks we need elsewhere).
# ---------------------------------------------------------------------

# G8 (was the red "DELAY" cell) -> becomes the orange "Tạm ổn" look
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)

# G7, G9, G10, G13 -> become the green "ok" look (same look G11 already has)
$ws.Range("G11").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("G11").Copy()
$ws.Range("G9").PasteSpecial(-4122)

$ws.Range("G11").Copy()
$ws.Range("G10").PasteSpecial(-4122)

$ws.Range("G11").Copy()
$ws.Range("G13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Phase 2: introduce the brand-new text values first, in the precise
# order they should land in the shared-string table.
# ---------------------------------------------------------------------

$ws.Range("G12").Value = "nearly done"
$ws.Range("C13").Value = "Hưng"
$ws.Range("C14").Value = "Hoang & Hung"
$ws.Range("B13").Value = "Module liệt kê bộ ba"
$ws.Range("B14").Value = "Module nhận biết thực thể"
$ws.Range("B15").Value = "Module sinh câu truy vấn"
$ws.Range("D14").Value = "12/12"

# ---------------------------------------------------------------------
# Phase 3: remaining value updates (re-using already interned strings).
# ---------------------------------------------------------------------

$ws.Range("G7").Value = "ok"
$ws.Range("G8").Value = "Tạm ổn"
$ws.Range("G9").Value = "ok"
$ws.Range("G10").Value = "ok"
$ws.Range("G13").Value = "ok"

$ws.Range("D13").Value = "26/11"
$ws.Range("E13").Value = 0.9

$ws.Range("C15").Value = "Hoang & Hung"
$ws.Range("D15").Value = "12/12"

$ws.Range("E14").Value = 0
$ws.Range("E15").Value = 0

$ws.Range("E12").Value = 0.9

$ws.Range("G14").Value = "on processing"
$ws.Range("G15").Value = "on processing"

# ---------------------------------------------------------------------
# Phase 4: selection / view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------

$ws.Select()
$ws.Range("G13").Select()
